# Update the "Final Value" column (D) on the active sheet with newly
# computed results, as part of building the polar graph of results (#75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2119.75
$ws.Range("D3").Value = 62.98161450745967
$ws.Range("D4").Value = 13.78116575408985
$ws.Range("D5").Value = 6.383432963279248
$ws.Range("D6").Value = 32.1
$ws.Range("D7").Value = 178.25
$ws.Range("D8").Value = 888.66
$ws.Range("D9").Value = 0.2461445722610137
$ws.Range("D10").Value = 2.982200565965071
$ws.Range("D11").Value = 0.8317844644262485
$ws.Range("D12").Value = 2857.59
$ws.Range("D13").Value = 3510.83
$ws.Range("D14").Value = 15.4
$ws.Range("D15").Value = 305.16
$ws.Range("D16").Value = 1.57
$ws.Range("D17").Value = 1.96
$ws.Range("D18").Value = 1.31
$ws.Range("D19").Value = 15.94
$ws.Range("D20").Value = 3.74
$ws.Range("D21").Value = 66.2
$ws.Range("D22").Value = 1057.74
$ws.Range("D23").Value = 1.29
$ws.Range("D24").Value = 667.76
$ws.Range("D25").Value = 1180.86
$ws.Range("D26").Value = 10.75
